# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 379
$wsExpo.Range("F3").Value = 2105
$wsExpo.Range("F4").Value = 109

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 379
$wsAll.Range("F7").Value = 2105
$wsAll.Range("F8").Value = 109
